$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 16:12"

# Reorder countries (Sudafrica, Mauricio, Consejo Danes para los Refugiados moved)
# and refresh updated case statistics
$ws.Cells.Item(8,2).Value = 40585
$ws.Cells.Item(8,3).Value = 3262
$ws.Cells.Item(8,4).Value = 5669
$ws.Cells.Item(8,5).Value = 34687
$ws.Cells.Item(33,5).Value = 1100
$ws.Cells.Item(33,7).Value = 1
$ws.Cells.Item(33,8).Value = 9
$ws.Cells.Item(39,1).Value = "Sudafrica"
$ws.Cells.Item(39,2).Value = 900
$ws.Cells.Item(39,3).Value = 191
$ws.Cells.Item(39,4).Value = 12
$ws.Cells.Item(39,5).Value = 888
$ws.Cells.Item(39,6).Value = 2
$ws.Cells.Item(39,7).Value = 0
$ws.Cells.Item(39,8).Value = 0
$ws.Cells.Item(40,1).Value = "Indonesia"
$ws.Cells.Item(40,2).Value = 893
$ws.Cells.Item(40,3).Value = 103
$ws.Cells.Item(40,4).Value = 35
$ws.Cells.Item(40,5).Value = 780
$ws.Cells.Item(40,6).Value = 0
$ws.Cells.Item(40,7).Value = 20
$ws.Cells.Item(40,8).Value = 78
$ws.Cells.Item(41,1).Value = "Rusia"
$ws.Cells.Item(41,2).Value = 840
$ws.Cells.Item(41,3).Value = 182
$ws.Cells.Item(41,4).Value = 38
$ws.Cells.Item(41,5).Value = 799
$ws.Cells.Item(41,6).Value = 8
$ws.Cells.Item(41,7).Value = 0
$ws.Cells.Item(41,8).Value = 3
$ws.Cells.Item(42,1).Value = "Grecia"
$ws.Cells.Item(42,2).Value = 821
$ws.Cells.Item(42,3).Value = 0
$ws.Cells.Item(42,4).Value = 36
$ws.Cells.Item(42,5).Value = 762
$ws.Cells.Item(42,6).Value = 53
$ws.Cells.Item(42,7).Value = 1
$ws.Cells.Item(42,8).Value = 23
$ws.Cells.Item(43,1).Value = "Islandia"
$ws.Cells.Item(43,2).Value = 802
$ws.Cells.Item(43,3).Value = 65
$ws.Cells.Item(43,4).Value = 68
$ws.Cells.Item(43,5).Value = 732
$ws.Cells.Item(43,6).Value = 11
$ws.Cells.Item(43,7).Value = 0
$ws.Cells.Item(43,8).Value = 2
$ws.Cells.Item(44,1).Value = "India"
$ws.Cells.Item(44,2).Value = 719
$ws.Cells.Item(44,3).Value = 62
$ws.Cells.Item(44,4).Value = 45
$ws.Cells.Item(44,5).Value = 658
$ws.Cells.Item(44,6).Value = 0
$ws.Cells.Item(44,7).Value = 4
$ws.Cells.Item(44,8).Value = 16
$ws.Cells.Item(45,1).Value = "Crucero"
$ws.Cells.Item(45,2).Value = 712
$ws.Cells.Item(45,4).Value = 597
$ws.Cells.Item(45,5).Value = 105
$ws.Cells.Item(45,6).Value = 15
$ws.Cells.Item(45,8).Value = 10
$ws.Cells.Item(51,5).Value = 512
$ws.Cells.Item(51,7).Value = 1
$ws.Cells.Item(51,8).Value = 6
$ws.Cells.Item(96,2).Value = 111
$ws.Cells.Item(96,3).Value = 30
$ws.Cells.Item(96,5).Value = 108
$ws.Cells.Item(96,7).Value = 1
$ws.Cells.Item(96,8).Value = 1
$ws.Cells.Item(105,1).Value = "Mauricio"
$ws.Cells.Item(105,2).Value = 81
$ws.Cells.Item(105,3).Value = 33
$ws.Cells.Item(105,4).Value = 0
$ws.Cells.Item(105,5).Value = 79
$ws.Cells.Item(105,6).Value = 1
$ws.Cells.Item(105,8).Value = 2
$ws.Cells.Item(106,1).Value = "Costa de Marfil"
$ws.Cells.Item(106,2).Value = 80
$ws.Cells.Item(106,3).Value = 0
$ws.Cells.Item(106,4).Value = 3
$ws.Cells.Item(106,5).Value = 77
$ws.Cells.Item(106,6).Value = 0
$ws.Cells.Item(107,1).Value = "Georgia"
$ws.Cells.Item(107,2).Value = 79
$ws.Cells.Item(107,3).Value = 4
$ws.Cells.Item(107,4).Value = 10
$ws.Cells.Item(107,5).Value = 69
$ws.Cells.Item(107,6).Value = 1
$ws.Cells.Item(107,8).Value = 0
$ws.Cells.Item(108,1).Value = "Camerun"
$ws.Cells.Item(108,2).Value = 75
$ws.Cells.Item(108,4).Value = 2
$ws.Cells.Item(108,6).Value = 0
$ws.Cells.Item(109,1).Value = "Guadalupe"
$ws.Cells.Item(109,2).Value = 73
$ws.Cells.Item(109,3).Value = 0
$ws.Cells.Item(109,5).Value = 72
$ws.Cells.Item(109,6).Value = 4
$ws.Cells.Item(110,1).Value = "Montenegro"
$ws.Cells.Item(110,2).Value = 67
$ws.Cells.Item(110,3).Value = 14
$ws.Cells.Item(110,5).Value = 66
$ws.Cells.Item(110,6).Value = 1
$ws.Cells.Item(111,1).Value = "Martinica"
$ws.Cells.Item(111,2).Value = 66
$ws.Cells.Item(111,3).Value = 0
$ws.Cells.Item(111,6).Value = 7
$ws.Cells.Item(111,8).Value = 1
$ws.Cells.Item(112,1).Value = "Uzbekistan"
$ws.Cells.Item(112,2).Value = 65
$ws.Cells.Item(112,3).Value = 5
$ws.Cells.Item(112,5).Value = 65
$ws.Cells.Item(112,6).Value = 4
$ws.Cells.Item(112,8).Value = 0
$ws.Cells.Item(113,1).Value = "Trinidad yTobago"
$ws.Cells.Item(113,2).Value = 61
$ws.Cells.Item(113,3).Value = 1
$ws.Cells.Item(113,4).Value = 0
$ws.Cells.Item(113,5).Value = 60
$ws.Cells.Item(113,6).Value = 0
$ws.Cells.Item(114,1).Value = "Cuba"
$ws.Cells.Item(114,2).Value = 57
$ws.Cells.Item(114,4).Value = 1
$ws.Cells.Item(114,5).Value = 55
$ws.Cells.Item(114,6).Value = 2
$ws.Cells.Item(114,8).Value = 1
$ws.Cells.Item(115,1).Value = "Honduras"
$ws.Cells.Item(115,3).Value = 0
$ws.Cells.Item(115,5).Value = 52
$ws.Cells.Item(115,6).Value = 0
$ws.Cells.Item(115,8).Value = 0
$ws.Cells.Item(117,1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(117,3).Value = 3
$ws.Cells.Item(117,4).Value = 0
$ws.Cells.Item(117,7).Value = 1
$ws.Cells.Item(117,8).Value = 3
$ws.Cells.Item(118,1).Value = "Nigeria"
$ws.Cells.Item(118,3).Value = 0
$ws.Cells.Item(118,4).Value = 2
$ws.Cells.Item(118,7).Value = 0
$ws.Cells.Item(118,8).Value = 1
